$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IDA")
Write-Host ("Col A width: " + $ws.Columns.Item(1).ColumnWidth())
Write-Host ("Col C width before: " + $ws.Columns.Item(3).ColumnWidth())
Write-Host ("Col D width before: " + $ws.Columns.Item(4).ColumnWidth())
Write-Host ("Col G width before: " + $ws.Columns.Item(7).ColumnWidth())
